# Auto-generated script applying Halicarnassus_Profits value updates
$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value2 = 2321.9375
$ws.Range("I38").Value2 = 65.09999999999999
$ws.Range("K38").Value2 = 195.3
$ws.Range("M38").Value2 = 176.7
# Row 64
$ws.Range("H64").Value2 = 0
$ws.Range("J64").Value2 = 0
$ws.Range("L64").Value2 = 0
$ws.Range("N64").ClearContents()
# Row 67
$ws.Range("H67").Value2 = 0
$ws.Range("J67").Value2 = 0
$ws.Range("L67").Value2 = 0
$ws.Range("N67").ClearContents()
# Row 74
$ws.Range("H74").Value2 = 11632.889
$ws.Range("I74").Value2 = 6528
$ws.Range("K74").Value2 = 6528
$ws.Range("M74").Value2 = -5592
# Row 77
$ws.Range("H77").Value2 = 11632.889
$ws.Range("I77").Value2 = 6528
$ws.Range("K77").Value2 = 32640
$ws.Range("M77").Value2 = -27960
# Row 105
$ws.Range("H105").Value2 = 15531
$ws.Range("J105").Value2 = 15531
$ws.Range("L105").Value2 = 15531
$ws.Range("N105").Value2 = -22519

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value2 = 1065.3334
$ws.Range("I2").Value2 = 1053
$ws.Range("J2").Value2 = 1090
$ws.Range("K2").Value2 = 1053
$ws.Range("L2").Value2 = 1090
$ws.Range("M2").Value2 = -940
$ws.Range("N2").Value2 = -1316
# Row 19
$ws.Range("H19").Value2 = 8125
$ws.Range("I19").Value2 = 5156.25
$ws.Range("K19").Value2 = 5156.25
$ws.Range("M19").Value2 = -4927.25
# Row 34
$ws.Range("H34").Value2 = 250
$ws.Range("I34").Value2 = 250
$ws.Range("J34").Value2 = 0
$ws.Range("K34").Value2 = 250
$ws.Range("L34").Value2 = 0
$ws.Range("M34").Value2 = 21
$ws.Range("N34").ClearContents()
# Row 61
$ws.Range("H61").Value2 = 2842.7144
$ws.Range("I61").Value2 = 1781
$ws.Range("K61").Value2 = 1781
$ws.Range("M61").Value2 = -1569
# Row 88
$ws.Range("H88").Value2 = 1306.5385
$ws.Range("J88").Value2 = 967.25
$ws.Range("L88").Value2 = 967.25
$ws.Range("N88").Value2 = -1779.25
# Row 91
$ws.Range("H91").Value2 = 1306.5385
$ws.Range("J91").Value2 = 967.25
$ws.Range("L91").Value2 = 967.25
$ws.Range("N91").Value2 = -3775.25
# Row 116
$ws.Range("H116").Value2 = 1065.3334
$ws.Range("I116").Value2 = 1053
$ws.Range("J116").Value2 = 1090
$ws.Range("K116").Value2 = 1053
$ws.Range("L116").Value2 = 1090
$ws.Range("M116").Value2 = 1241
$ws.Range("N116").Value2 = -5678
# Row 122
$ws.Range("H122").Value2 = 1577.3334
$ws.Range("I122").Value2 = 1293
$ws.Range("K122").Value2 = 3879
$ws.Range("M122").Value2 = -1429
# Row 132
$ws.Range("H132").Value2 = 856.37036
$ws.Range("J132").Value2 = 831
$ws.Range("L132").Value2 = 2493
$ws.Range("N132").Value2 = -7553
# Row 136
$ws.Range("H136").Value2 = 2842.7144
$ws.Range("I136").Value2 = 1781
$ws.Range("K136").Value2 = 5343
$ws.Range("M136").Value2 = -2793

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value2 = 1065.3334
$ws.Range("I3").Value2 = 1053
$ws.Range("J3").Value2 = 1090
$ws.Range("K3").Value2 = 1053
$ws.Range("L3").Value2 = 1090
$ws.Range("M3").Value2 = -939
$ws.Range("N3").Value2 = -1318
# Row 11
$ws.Range("H11").Value2 = 451.66666
$ws.Range("I11").Value2 = 400
$ws.Range("J11").Value2 = 477.5
$ws.Range("K11").Value2 = 400
$ws.Range("L11").Value2 = 477.5
$ws.Range("M11").Value2 = -260
$ws.Range("N11").Value2 = -757.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 39
$ws.Range("H39").Value2 = 2535.5
$ws.Range("I39").Value2 = 2535.5
$ws.Range("K39").Value2 = 2535.5
$ws.Range("M39").Value2 = -2144.5
# Row 49
$ws.Range("H49").Value2 = 2535.5
$ws.Range("I49").Value2 = 2535.5
$ws.Range("K49").Value2 = 2535.5
$ws.Range("M49").Value2 = -2353.5
# Row 99
$ws.Range("H99").Value2 = 1735.2142
$ws.Range("I99").Value2 = 1616.1666
$ws.Range("J99").Value2 = 2449.5
$ws.Range("K99").Value2 = 1616.1666
$ws.Range("L99").Value2 = 2449.5
$ws.Range("M99").Value2 = -118.1666
$ws.Range("N99").Value2 = -5445.5
# Row 107
$ws.Range("H107").Value2 = 1833.909
$ws.Range("I107").Value2 = 625.7143
$ws.Range("J107").Value2 = 3948.25
$ws.Range("K107").Value2 = 625.7143
$ws.Range("L107").Value2 = 3948.25
$ws.Range("M107").Value2 = 1294.2857
$ws.Range("N107").Value2 = -7788.25
# Row 122
$ws.Range("H122").Value2 = 1803.2
$ws.Range("I122").Value2 = 1782.4445
$ws.Range("K122").Value2 = 5347.333500000001
$ws.Range("M122").Value2 = -2897.333500000001
# Row 126
$ws.Range("H126").Value2 = 1735.2142
$ws.Range("I126").Value2 = 1616.1666
$ws.Range("J126").Value2 = 2449.5
$ws.Range("K126").Value2 = 4848.4998
$ws.Range("L126").Value2 = 7348.5
$ws.Range("M126").Value2 = -2378.4998
$ws.Range("N126").Value2 = -12288.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value2 = 1907.0769
$ws.Range("I5").Value2 = 1473.25
$ws.Range("K5").Value2 = 4419.75
$ws.Range("M5").Value2 = -4307.75
# Row 49
$ws.Range("H49").Value2 = 0
$ws.Range("I49").Value2 = 0
$ws.Range("K49").Value2 = 0
$ws.Range("M49").ClearContents()
# Row 59
$ws.Range("H59").Value2 = 900
$ws.Range("I59").Value2 = 900
$ws.Range("K59").Value2 = 2700
$ws.Range("M59").Value2 = -2160
# Row 60
$ws.Range("H60").Value2 = 1494.2307
$ws.Range("J60").Value2 = 2271.875
$ws.Range("L60").Value2 = 6815.625
$ws.Range("N60").Value2 = -7317.625
# Row 68
$ws.Range("H68").Value2 = 439.6
$ws.Range("I68").Value2 = 399.33334
$ws.Range("K68").Value2 = 1198.00002
$ws.Range("M68").Value2 = -387.0000199999999
# Row 71
$ws.Range("H71").Value2 = 439.6
$ws.Range("I71").Value2 = 399.33334
$ws.Range("K71").Value2 = 3594.00006
$ws.Range("M71").Value2 = 461.9999399999997
# Row 80
$ws.Range("H80").Value2 = 4335.4346
$ws.Range("I80").Value2 = 4021.889
$ws.Range("J80").Value2 = 5464.2
$ws.Range("K80").Value2 = 12065.667
$ws.Range("L80").Value2 = 16392.6
$ws.Range("M80").Value2 = -11129.667
$ws.Range("N80").Value2 = -18264.6
# Row 83
$ws.Range("H83").Value2 = 4335.4346
$ws.Range("I83").Value2 = 4021.889
$ws.Range("J83").Value2 = 5464.2
$ws.Range("K83").Value2 = 36197.001
$ws.Range("L83").Value2 = 49177.8
$ws.Range("M83").Value2 = -31517.001
$ws.Range("N83").Value2 = -58537.8
# Row 129
$ws.Range("H129").Value2 = 1518.8572
$ws.Range("J129").Value2 = 1544
$ws.Range("L129").Value2 = 4632
$ws.Range("N129").Value2 = -14632
# Row 132
$ws.Range("H132").Value2 = 3440.1538
$ws.Range("I132").Value2 = 2684.125
$ws.Range("K132").Value2 = 24157.125
$ws.Range("M132").Value2 = -21627.125
# Row 135
$ws.Range("H135").Value2 = 1907.0769
$ws.Range("I135").Value2 = 1473.25
$ws.Range("K135").Value2 = 13259.25
$ws.Range("M135").Value2 = -10724.25

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 58
$ws.Range("H58").Value2 = 24000
$ws.Range("I58").Value2 = 0
$ws.Range("J58").Value2 = 24000
$ws.Range("K58").Value2 = 0
$ws.Range("L58").Value2 = 24000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value2 = -24554
# Row 113
$ws.Range("H113").Value2 = 3353.2222
$ws.Range("I113").Value2 = 2282.7144
$ws.Range("J113").Value2 = 7100
$ws.Range("K113").Value2 = 2282.7144
$ws.Range("L113").Value2 = 7100
$ws.Range("M113").Value2 = -112.7143999999998
$ws.Range("N113").Value2 = -11440
# Row 122
$ws.Range("H122").Value2 = 2594.4
$ws.Range("I122").Value2 = 2394
$ws.Range("J122").Value2 = 2644.5
$ws.Range("K122").Value2 = 7182
$ws.Range("L122").Value2 = 7933.5
$ws.Range("M122").Value2 = -4732
$ws.Range("N122").Value2 = -12833.5
# Row 126
$ws.Range("H126").Value2 = 2362.25
$ws.Range("I126").Value2 = 1474.75
$ws.Range("K126").Value2 = 4424.25
$ws.Range("M126").Value2 = -1954.25

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 17
$ws.Range("H17").Value2 = 0
$ws.Range("I17").Value2 = 0
$ws.Range("K17").Value2 = 0
$ws.Range("M17").ClearContents()
# Row 26
$ws.Range("H26").Value2 = 3500
$ws.Range("I26").Value2 = 3500
$ws.Range("K26").Value2 = 3500
$ws.Range("M26").Value2 = -3205
# Row 46
$ws.Range("H46").Value2 = 8566.666999999999
$ws.Range("I46").Value2 = 1400
$ws.Range("J46").Value2 = 10000
$ws.Range("K46").Value2 = 1400
$ws.Range("L46").Value2 = 10000
$ws.Range("M46").Value2 = -1212
$ws.Range("N46").Value2 = -10376
# Row 50
$ws.Range("H50").Value2 = 66078
$ws.Range("I50").Value2 = 66078
$ws.Range("K50").Value2 = 66078
$ws.Range("M50").Value2 = -65441
# Row 54
$ws.Range("H54").Value2 = 20000
$ws.Range("J54").Value2 = 20000
$ws.Range("L54").Value2 = 20000
$ws.Range("N54").Value2 = -21288
# Row 106
$ws.Range("H106").Value2 = 7478.1665
$ws.Range("J106").Value2 = 7478.1665
$ws.Range("L106").Value2 = 7478.1665
$ws.Range("N106").Value2 = -10002.1665
# Row 128
$ws.Range("H128").Value2 = 40000
$ws.Range("J128").Value2 = 40000
$ws.Range("L128").Value2 = 40000
$ws.Range("N128").Value2 = -49960

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value2 = 443.33334
$ws.Range("I7").Value2 = 750
$ws.Range("J7").Value2 = 290
$ws.Range("K7").Value2 = 750
$ws.Range("L7").Value2 = 290
$ws.Range("M7").Value2 = -637
$ws.Range("N7").Value2 = -516
# Row 26
$ws.Range("H26").Value2 = 512
$ws.Range("I26").Value2 = 512
$ws.Range("J26").Value2 = 0
$ws.Range("K26").Value2 = 512
$ws.Range("L26").Value2 = 0
$ws.Range("M26").Value2 = -219
$ws.Range("N26").ClearContents()
# Row 37
$ws.Range("H37").Value2 = 20000
$ws.Range("J37").Value2 = 20000
$ws.Range("L37").Value2 = 20000
$ws.Range("N37").Value2 = -20406
# Row 100
$ws.Range("H100").Value2 = 1469.6
$ws.Range("I100").Value2 = 1410.6666
$ws.Range("K100").Value2 = 2821.3332
$ws.Range("M100").Value2 = -2280.3332
# Row 113
$ws.Range("H113").Value2 = 641.5
$ws.Range("I113").Value2 = 470.85715
$ws.Range("K113").Value2 = 1412.57145
$ws.Range("M113").Value2 = 757.4285500000001
# Row 122
$ws.Range("H122").Value2 = 2061.0833
$ws.Range("J122").Value2 = 2463.75
$ws.Range("L122").Value2 = 7391.25
$ws.Range("N122").Value2 = -12291.25
